$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = '69.075.19'
$ws.Range("E2").Value2 = '  -2.19%  '
$ws.Range("D3").Value2 = '3.678.18'
$ws.Range("D4").Value2 = '''0.999'
$ws.Range("E4").Value2 = '  -0.07%  '
$ws.Range("D5").Value2 = '''682.00'
$ws.Range("E5").Value2 = '  -3.52%  '
$ws.Range("D6").Value2 = '''162.49'
$ws.Range("E6").Value2 = '  -4.34%  '
$ws.Range("D7").Value2 = '3.674.48'
$ws.Range("E7").Value2 = '  -3.06%  '
$ws.Range("E8").Value2 = '  +0.00%  '
$ws.Range("D9").Value2 = '''0.496'
$ws.Range("E9").Value2 = '  -4.69%  '
$ws.Range("E10").Value2 = '  -7.49%  '
$ws.Range("D11").Value2 = '''7.24'
$ws.Range("E11").Value2 = '  -1.67%  '
$ws.Range("E12").Value2 = '  -1.39%  '
$ws.Range("E13").Value2 = '  -7.21%  '
$ws.Range("D14").Value2 = '''33.41'
$ws.Range("E14").Value2 = '  -7.70%  '
$ws.Range("D15").Value2 = '4.297.53'
$ws.Range("E15").Value2 = '  -3.00%  '
$ws.Range("D16").Value2 = '3.679.19'
$ws.Range("E16").Value2 = '  -3.41%  '
$ws.Range("D17").Value2 = '69.162.88'
$ws.Range("E17").Value2 = '  -2.06%  '
$ws.Range("E18").Value2 = '  -1.75%  '
$ws.Range("D19").Value2 = '''16.32'
$ws.Range("E19").Value2 = '  -5.98%  '
$ws.Range("D20").Value2 = '''6.63'
$ws.Range("E20").Value2 = '  -7.15%  '
$ws.Range("D21").Value2 = '''482.74'
$ws.Range("E21").Value2 = '  -2.17%  '
$ws.Range("D22").Value2 = '''9.81'
$ws.Range("E22").Value2 = '  -7.64%  '
$ws.Range("D23").Value2 = '''0.664'
$ws.Range("E23").Value2 = '  -8.61%  '
$ws.Range("D24").Value2 = '''79.50'
$ws.Range("E24").Value2 = '  -6.27%  '
$ws.Range("D25").Value2 = '3.820.97'
$ws.Range("E25").Value2 = '  -3.06%  '
$ws.Range("D26").Value2 = '''11.55'
$ws.Range("E26").Value2 = '  -4.33%  '
$ws.Range("B27").Value2 = 'PEPE'
$ws.Range("C27").Value2 = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D27").Value2 = '''0.0000127'
$ws.Range("E27").Value2 = '  -11.92%  '
$ws.Range("B28").Value2 = 'Dai'
$ws.Range("C28").Value2 = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D28").Value2 = '''0.999'
$ws.Range("E28").Value2 = '  -0.04%  '
$ws.Range("D29").Value2 = '''9.54'
$ws.Range("E29").Value2 = '  -8.75%  '
$ws.Range("D30").Value2 = '''1.84'
$ws.Range("E30").Value2 = '  -10.12%  '
$ws.Range("E31").Value2 = '  -11.47%  '
$ws.Range("D32").Value2 = '''2.11'
$ws.Range("E32").Value2 = '  -4.80%  '
$ws.Range("D33").Value2 = '''6.74'
$ws.Range("E33").Value2 = '  -7.98%  '
$ws.Range("E34").Value2 = '  +0.08%  '
$ws.Range("D35").Value2 = '''26.90'
$ws.Range("E35").Value2 = '  -7.52%  '
$ws.Range("E36").Value2 = '  -6.95%  '
$ws.Range("D37").Value2 = '3.640.88'
$ws.Range("E37").Value2 = '  -3.15%  '
$ws.Range("D38").Value2 = '''8.51'
$ws.Range("E38").Value2 = '  -5.93%  '
$ws.Range("E39").Value2 = '  +2.05%  '
$ws.Range("D40").Value2 = '''0.0944'
$ws.Range("E40").Value2 = '  -6.78%  '
$ws.Range("E43").Value2 = '  +0.12%  '
$ws.Range("D44").Value2 = '''0.956'
$ws.Range("E44").Value2 = '  -7.98%  '
$ws.Range("D45").Value2 = '''156.51'
$ws.Range("E45").Value2 = '  -4.91%  '
$ws.Range("E46").Value2 = '  -1.86%  '
$ws.Range("E47").Value2 = '  -14.72%  '
$ws.Range("B48").Value2 = 'FLOKI'
$ws.Range("C48").Value2 = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D48").Value2 = '''0.000277'
$ws.Range("E48").Value2 = '  -11.25%  '
$ws.Range("B49").Value2 = 'ONDO'
$ws.Range("C49").Value2 = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D49").Value2 = '''1.30'
$ws.Range("E49").Value2 = '  -4.45%  '
$ws.Range("B50").Value2 = 'Bittensor'
$ws.Range("C50").Value2 = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D50").Value2 = '''389.85'
$ws.Range("E50").Value2 = '  -7.41%  '
$ws.Range("D51").Value2 = '''28.12'
$ws.Range("E51").Value2 = '  +1.28%  '
